# Updates match-by-match batting activity figures (runs/balls/fours/sixes)
# for Shreyas Iyer (c) on the "Shreyas Iyer (c)" sheet, per the commit
# "updated activity till excel form". Values are written with a leading
# apostrophe so they stay text cells (matching the source data's stored
# type) rather than being auto-converted to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = "'7"
$ws.Range("D2").Value = "'9"
$ws.Range("E2").Value = "'0"
$ws.Range("F2").Value = "'0"

# Row 3
$ws.Range("C3").Value = "'25"
$ws.Range("D3").Value = "'29"
$ws.Range("E3").Value = "'1"
$ws.Range("F3").Value = "'1"

# Row 4
$ws.Range("D4").Value = "'12"

# Row 6
$ws.Range("C6").Value = "'12"
$ws.Range("D6").Value = "'8"
$ws.Range("E6").Value = "'3"
$ws.Range("F6").Value = "'0"

# Row 8
$ws.Range("C8").Value = "'53"
$ws.Range("D8").Value = "'43"
$ws.Range("E8").Value = "'3"
$ws.Range("F8").Value = "'2"

# Row 9
$ws.Range("C9").Value = "'23"
$ws.Range("D9").Value = "'23"

# Row 10
$ws.Range("C10").Value = "'42"
$ws.Range("D10").Value = "'33"
$ws.Range("E10").Value = "'5"

# Row 11
$ws.Range("C11").Value = "'14"
$ws.Range("D11").Value = "'12"
$ws.Range("E11").Value = "'0"
$ws.Range("F11").Value = "'1"

# Row 12
$ws.Range("C12").Value = "'65"
$ws.Range("D12").Value = "'50"
$ws.Range("E12").Value = "'6"
$ws.Range("F12").Value = "'2"

# Row 13
$ws.Range("C13").Value = "'17"
$ws.Range("D13").Value = "'21"
$ws.Range("E13").Value = "'2"

# Row 14
$ws.Range("C14").Value = "'88"
$ws.Range("D14").Value = "'38"
$ws.Range("E14").Value = "'7"
$ws.Range("F14").Value = "'6"

# Row 15
$ws.Range("C15").Value = "'26"
$ws.Range("D15").Value = "'22"

# Row 16
$ws.Range("C16").Value = "'39"
$ws.Range("D16").Value = "'32"
$ws.Range("E16").Value = "'0"
$ws.Range("F16").Value = "'3"

# Row 17
$ws.Range("C17").Value = "'11"
$ws.Range("D17").Value = "'13"
$ws.Range("E17").Value = "'1"

# Row 18
$ws.Range("C18").Value = "'22"
$ws.Range("D18").Value = "'18"
$ws.Range("E18").Value = "'4"
$ws.Range("F18").Value = "'0"
